$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("residual_conn")

# Header row for "Batch normalization" section (row 16)
$ws.Range("A16").Value = "Batch normalization"
$ws.Range("D18").Select()

# Data rows 17-18
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 3
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 0.96243876218795699

$ws.Range("A18").Value = 64
$ws.Range("B18").Value = 4
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 0.96316456794738703
